# Update gh-pages output (苏州-漫展信息.xlsx) to match the newly scraped data:
#  - refresh "want to go" (F column) counts across sheets 1 (展览) and 4 (全部类型)
#  - add the new "苏州·英雄时代2024哈瓦西钢琴演奏会" event row to sheet 2 (演出)
#    and to sheet 4 (全部类型), shifting the later rows of sheet 4 down by one

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition) - refreshed "want to go" counts only
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$sheet1FChanges = @{
    2 = 598; 3 = 495; 4 = 1275; 5 = 1131; 6 = 14182; 7 = 15929; 8 = 13; 9 = 67;
    17 = 35; 18 = 91; 19 = 31; 20 = 1230; 21 = 132; 23 = 24; 24 = 6323; 26 = 1105;
    27 = 5623; 29 = 144; 30 = 135; 31 = 4610; 32 = 8
}
foreach ($row in $sheet1FChanges.Keys) {
    $ws1.Range("F$row").Value = $sheet1FChanges[$row]
}

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performance) - append new row 3 for the new event
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Match the index column's bold/bordered style used by A2
$ws2.Range("A2").Copy()
$ws2.Range("A3").PasteSpecial(-4122)
$ws2.Range("A3").Value = 2

# The date column stores plain text ("2024-06-02"); force text format first so
# Excel does not silently convert it to a date serial number.
$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "2024-06-02"

$ws2.Range("C3").Value = "苏州·英雄时代2024哈瓦西钢琴演奏会"
$ws2.Range("D3").Value = "东太湖大道12000号 苏州湾大剧院"
$ws2.Range("E3").Value = "2024.06.02 19:30-06.02 21:00"
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 499
$ws2.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=83901"
$ws2.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202404/LbCirky11712569675168.png"

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) - refresh counts, then insert the new event
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$sheet4FChanges = @{
    2 = 598; 3 = 495; 4 = 1275; 5 = 1131; 6 = 14182; 7 = 15929; 8 = 13; 9 = 67;
    17 = 35; 18 = 91; 19 = 31; 20 = 1230; 21 = 132; 24 = 24; 25 = 6323; 27 = 1105
}
foreach ($row in $sheet4FChanges.Keys) {
    $ws4.Range("F$row").Value = $sheet4FChanges[$row]
}

# Insert a new blank row at position 28 (pushes old rows 28-33 down to 29-34)
$ws4.Rows("28:28").Insert()

# The inserted row's index cell lost its style; copy it back from the row above
$ws4.Range("A27").Copy()
$ws4.Range("A28").PasteSpecial(-4122)
$ws4.Range("A28").Value = 27

$ws4.Range("B28").NumberFormat = "@"
$ws4.Range("B28").Value = "2024-06-02"

$ws4.Range("C28").Value = "苏州·英雄时代2024哈瓦西钢琴演奏会"
$ws4.Range("D28").Value = "东太湖大道12000号 苏州湾大剧院"
$ws4.Range("E28").Value = "2024.06.02 19:30-06.02 21:00"
$ws4.Range("F28").Value = 0
$ws4.Range("G28").Value = 499
$ws4.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=83901"
$ws4.Range("I28").Value = "//i0.hdslb.com/bfs/openplatform/202404/LbCirky11712569675168.png"

# Refreshed "want to go" counts for the rows that shifted down one position
$ws4.Range("F29").Value = 5623
$ws4.Range("F31").Value = 144
$ws4.Range("F32").Value = 135
$ws4.Range("F33").Value = 4611
$ws4.Range("F34").Value = 8

Write-Host "Edit complete"
